$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Done list (numId=2): the old "Documentation" bullet is dropped, and a new
# bullet "Bugfix throw bottle: Bottle follows character after throw" is
# appended after "Collision".
# (Do this before touching the Todo list, so the later "Documentation"
# search below unambiguously finds the Done-list occurrence.)
# ------------------------------------------------------------------
$rOldDoc = $d.Content
$okOldDoc = $rOldDoc.Find.Execute("Documentation", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pOldDoc = $rOldDoc.Paragraphs(1)
$pOldDoc.Range.Delete()

$rCollision = $d.Content
$okCollision = $rCollision.Find.Execute("Collision", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pCollision = $rCollision.Paragraphs(1)

$pCollision.Range.InsertParagraphAfter()
$pBugfix = $pCollision.Next()
$pBugfix.Range.Text = "Bugfix throw bottle: Bottle follows character after throw"
$pBugfix.Range.ListFormat.ListLevelNumber = 1

# A new blank paragraph (no list formatting) is also added right after the
# bugfix bullet, before the document's pre-existing trailing blank paragraph.
$pBugfix.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# Todo list (numId=1): append two new bullet items after the last one,
# "Images in separate file" -> "Documentation", "Intro Screen".
# Inserting right after an existing list item makes the new paragraph
# inherit the Listenabsatz style plus the numId=1/ilvl=0 numbering.
# ------------------------------------------------------------------
$rImages = $d.Content
$okImages = $rImages.Find.Execute("Images in separate file", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pImages = $rImages.Paragraphs(1)

$pImages.Range.InsertParagraphAfter()
$pDocTodo = $pImages.Next()
$pDocTodo.Range.Text = "Documentation"

$pDocTodo.Range.InsertParagraphAfter()
$pIntro = $pDocTodo.Next()
$pIntro.Range.Text = "Intro Screen"

# Remove the old blank separator paragraph that used to sit between the Todo
# list and "Done" -- it is now the paragraph right after "Intro Screen".
$pBlank = $pIntro.Next()
$pBlank.Range.Delete()
